# Update the "ランサーズ" worksheet (rows 2-5 get new scrape results, rows 6-15 removed).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Row 2 ---
$ws.Range("A2").Value = "2025-11-08 06:24:04"
$ws.Range("B2").Value = "WordPress保守業務および保守業務の自動化ツール開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5429668"
$ws.Range("G2").Value = 218
$ws.Range("H2").Value = "◆ツール,開発 ○WordPress"

# --- Row 3 ---
$ws.Range("A3").Value = "2025-11-08 06:24:04"
$ws.Range("B3").Value = "複数の見積書から情報抜き出して別表に書き出す作業の自動化 Excel VBAツール開発依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5429304"
$ws.Range("G3").Value = 208
$ws.Range("H3").Value = "◆ツール,開発"

# --- Row 4 ---
$ws.Range("A4").Value = "2025-11-08 06:24:04"
$ws.Range("B4").Value = "【急募】既存で作成済みのAccessデータベースをPower Apps等のアプリに移行したい"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5429495"
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = "◇アプリ"

# --- Row 5 ---
$ws.Range("A5").Value = "2025-11-08 06:24:04"
$ws.Range("B5").Value = "【フルスタックエンジニア】 働きながらスキルアップもできるEC業界で活躍してみませんか?"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5429335"
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

# --- Remove the now-obsolete rows 6-15 ---
$ws.Rows("6:15").Delete()

# --- Hyperlinks: the engine's Hyperlinks.Delete() is collection-wide (not
# range-scoped), so clear everything and rebuild only the F2:F5 links that
# should survive, then restore the "Hyperlink" cell style on them.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5429668")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5429304")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5429495")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5429335")
$ws.Range("F2:F5").Style = "Hyperlink"

# --- Column width tweaks (col B 52->48, col D 30->28, col H 19->20). The
# engine's ColumnWidth setter bakes in Excel's character-width padding of
# +5/6 character before it is written back to the OOXML <col width=.../>,
# so subtract that offset to land on the exact target width.
$ws.Columns("B").ColumnWidth = 47.166666666666664
$ws.Columns("D").ColumnWidth = 27.166666666666668
$ws.Columns("H").ColumnWidth = 19.166666666666668
